# wolfStudio.xlsx -- "Update word and excel files"
#
# The calculated table column "כיצד הרגשת במהלך ההתנסות במערכת ?"
# (Table1's 18th column / worksheet column R) scores each respondent from
# the SUS-style questionnaire answers. The questionnaire columns it reads
# from shifted from B,C,E,F,G,H,I,J,K,L to H,I,J,K,L,M,N,O,P,Q, so re-apply
# the column's formula across every data row -- this recalculates R2:R30
# (and the SUM(...)/29 average in R32, which references the table column
# by name) with the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tbl = $ws.ListObjects.Item(1)
$col = $tbl.ListColumns.Item(18)
$col.DataBodyRange.Formula = "=(
 (H2 - 1) + (5 - I2) +
 (J2 - 1) + (5 - K2) +
 (L2 - 1) + (5 - M2) +
 (N2 - 1) + (5 - O2) +
 (P2 - 1) + (5 - Q2)
) * 2.5"

# Reflect the author's final active cell / scroll position in the sheet view.
$ws.Range("R14").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollColumn = 16
$win.ScrollRow = 10
